# alignment_manual_changes.xlsx — remove taxa with branch lengths > 0.9
#
# The "Family" sheet's Change/Comment columns record which taxon was
# dropped from each orthogroup alignment and why. Six orthogroups lost a
# taxon for exceeding the 0.9 branch-length cutoff:
#   - rows 8, 12, 13, 16, 20  -> LS484712_Daphniairidovirus_daphnia1
#   - row 32                  -> AY894343_Megalocytivirus_pagrus1_RSIV
# each with its own recorded branch length.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Family")

# Row 8 — OG000006
$ws.Range("C8").Value = "Branch length 0.9054"
$ws.Range("B8").Value = "Removed LS484712_Daphniairidovirus_daphnia1"

# Row 12 — G0000010
$ws.Range("B12").Value = "Removed LS484712_Daphniairidovirus_daphnia1"
$ws.Range("C12").Value = "Branch length 1.339"

# Row 13 — G0000011
$ws.Range("B13").Value = "Removed LS484712_Daphniairidovirus_daphnia1"
$ws.Range("C13").Value = "Branch length 0.9578"

# Row 16 — G0000014
$ws.Range("B16").Value = "Removed LS484712_Daphniairidovirus_daphnia1"
$ws.Range("C16").Value = "Branch length 1.181"

# Row 20 — G0000018
$ws.Range("B20").Value = "Removed LS484712_Daphniairidovirus_daphnia1"
$ws.Range("C20").Value = "Branch length 1.108"

# Row 32 — G0000030
$ws.Range("B32").Value = "Removed AY894343_Megalocytivirus_pagrus1_RSIV"
$ws.Range("C32").Value = "Branch length 1.049"

# Reflect the author's final on-screen state: Family sheet active/selected
# with B29 the active cell (Genus loses the tab-selected flag as a result).
[void]$ws.Activate()
[void]$ws.Range("B29").Select()
